$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "307.90"
    "E2" = "0.82%"
    "D3" = "36.27"
    "E3" = "1.14%"
    "D4" = "5.055"
    "E4" = "1.38%"
    "D5" = "0.08107"
    "E5" = "0.27%"
    "D6" = "2.020"
    "E6" = "5.97%"
    "D7" = "7.859"
    "E7" = "-0.23%"
    "D8" = "0.9271"
    "E8" = "-0.31%"
    "D9" = "0.1484"
    "E9" = "18.76%"
    "D10" = "0.1944"
    "E10" = "2.05%"
    "D11" = "0.09100"
    "E11" = "-1.25%"
    "D12" = "0.03522"
    "E12" = "-0.17%"
    "D13" = "0.09887"
    "E13" = "-0.42%"
    "D14" = "0.001417"
    "E14" = "0.24%"
    "D15" = "0.006062"
    "E15" = "0.17%"
    "D16" = "3.838"
    "E16" = "6.67%"
    "E17" = "0.55%"
    "D18" = "3.437"
    "E18" = "10.68%"
    "D19" = "0.3458"
    "E19" = "0.15%"
    "D20" = "0.1303"
    "E20" = "0.61%"
    "D21" = "4.823"
    "E21" = "-7.65%"
    "D22" = "0.2343"
    "E22" = "-7.47%"
    "D23" = "0.04397"
    "E23" = "-0.48%"
    "D24" = "0.001232"
    "E24" = "-0.23%"
    "D25" = "0.004189"
    "E27" = "0.10%"
    "D39" = "0.02058"
    "E39" = "5.33%"
    "D40" = "0.05139"
    "E40" = "-2.14%"
    "D41" = "0.007470"
    "E41" = "-0.94%"
    "D42" = "0.01005"
    "E42" = "-0.76%"
    "D43" = "0.1372"
    "E43" = "-0.05%"
    "D44" = "0.002121"
    "E44" = "1.06%"
    "D45" = "0.009891"
    "E45" = "-7.75%"
    "D46" = "0.00006301"
    "E46" = "-0.72%"
    "D47" = "0.00000000751"
    "E47" = "-0.01%"
    "D48" = "63.84"
    "E48" = "0.43%"
    "E49" = "-3.54%"
    "D50" = "0.00002102"
    "E50" = "-0.01%"
    "D51" = "0.0002002"
    "E51" = "-0.01%"
}

foreach ($addr in $updates.Keys) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$addr]
    $rng.Style = $origStyle
}

Write-Output ("Updated " + $updates.Count + " cells")